# The deck ships two DrawingML themes:
#   ppt/theme/theme1.xml -> used by the (only) slide master, currently the
#                            "Integral" palette, driving every slide.
#   ppt/theme/theme2.xml -> used by the notes master, currently the default
#                            "Office Theme" palette.
#
# The target commit swaps the two themes' contents: the slide master's
# theme (theme1.xml) becomes the stock "Office Theme" colour scheme, while
# the notes master's theme (theme2.xml) becomes the "Integral" colours.
# The font scheme and the fill/line/effect format scheme are byte-identical
# between the two themes already, so only the 12 colour-scheme slots
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) actually need to change.
#
# PowerPoint's automation surface doesn't expose a way to rename a theme /
# colour scheme in place, but it does let us repaint the twelve theme
# colour slots via ThemeColorScheme.Colors(i).RGB, which is exactly the
# substantive part of this edit (the <a:srgbClr val="..."/> values).

$p = $ppt.ActivePresentation

# VBA RGB(r,g,b) packs as r + g*256 + b*65536 ("0x00BBGGRR").
function RgbDec([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme palette (hex -> decimal), in ThemeColorScheme.Colors() order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
# 9 accent5, 10 accent6, 11 hlink, 12 folHlink.
$officeTheme = @(
    (RgbDec 0x00 0x00 0x00),  # dk1      000000
    (RgbDec 0xFF 0xFF 0xFF),  # lt1      FFFFFF
    (RgbDec 0x44 0x54 0x6A),  # dk2      44546A
    (RgbDec 0xE7 0xE6 0xE6),  # lt2      E7E6E6
    (RgbDec 0x5B 0x9B 0xD5),  # accent1  5B9BD5
    (RgbDec 0xED 0x7D 0x31),  # accent2  ED7D31
    (RgbDec 0xA5 0xA5 0xA5),  # accent3  A5A5A5
    (RgbDec 0xFF 0xC0 0x00),  # accent4  FFC000
    (RgbDec 0x44 0x72 0xC4),  # accent5  4472C4
    (RgbDec 0x70 0xAD 0x47),  # accent6  70AD47
    (RgbDec 0x05 0x63 0xC1),  # hlink    0563C1
    (RgbDec 0x95 0x4F 0x72)   # folHlink 954F72
)

# All slides share the single slide master -> theme1.xml, so repainting the
# scheme from any one slide updates the theme part used everywhere.
$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = $officeTheme[$i - 1]
}
